$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
$hdr = $ws.Range("A8")
$hdr.Characters(21,2).Text = "31"
$rpt = $ws.Range("C9")
$rpt.Characters(27,9).Text = "7/31/2023"
$rpt.Characters(47,9).Text = "8/6/2023"

# --- Cells changing from numeric to text (e.g. "0" or "***.*") ---
# Set the quoted text value first, then paste-format from a same-style donor cell
# so the final style index matches the target (avoids the quote-prefix style).
$ws.Range("F14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("D18").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F28").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F29").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G30").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("H30").PasteSpecial(-4122)

# --- Cells changing from text to numeric ---
$ws.Range("C18").Value = 4
$ws.Range("I14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C26").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = -33.333333333333
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 37.5
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = -15.702479338843
$ws.Range("L16").Value = 52.238805970149
$ws.Range("M16").Value = -37.423312883435
$ws.Range("N16").Value = -84.661654135338
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -7.407407407407
$ws.Range("I17").Value = 221
$ws.Range("J17").Value = 233
$ws.Range("K17").Value = -5.150214592274
$ws.Range("L17").Value = 40.764331210191
$ws.Range("M17").Value = 61.313868613138
$ws.Range("N17").Value = 6.763285024154
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 91
$ws.Range("K18").Value = 12.345679012345
$ws.Range("L18").Value = 31.884057971014
$ws.Range("M18").Value = -50.273224043715
$ws.Range("N18").Value = -88.996372430471
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -3.571428571428
$ws.Range("I19").Value = 187
$ws.Range("J19").Value = 214
$ws.Range("K19").Value = -12.616822429906
$ws.Range("L19").Value = 43.846153846153
$ws.Range("M19").Value = -3.608247422680
$ws.Range("N19").Value = -42.813455657492
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -62.5
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 162
$ws.Range("J20").Value = 152
$ws.Range("K20").Value = 6.578947368421
$ws.Range("L20").Value = 18.248175182481
$ws.Range("M20").Value = 2.531645569620
$ws.Range("N20").Value = -92.151162790697
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -4.347826086956
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = 3.669724770642
$ws.Range("I21").Value = 777
$ws.Range("J21").Value = 816
$ws.Range("K21").Value = -4.779411764705
$ws.Range("L21").Value = 34.429065743944
$ws.Range("M21").Value = -8.156028368794
$ws.Range("N21").Value = -81.108679795769
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -36.363636363636
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = -30
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -25.714285714285
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = -18.796992481203
$ws.Range("I24").Value = 893
$ws.Range("J24").Value = 913
$ws.Range("K24").Value = -2.190580503833
$ws.Range("L24").Value = 55.574912891986
$ws.Range("M24").Value = 90.811965811965
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 331
$ws.Range("J25").Value = 330
$ws.Range("K25").Value = 0.303030303030
$ws.Range("L25").Value = 18.637992831541
$ws.Range("M25").Value = -11.260053619302
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 22
$ws.Range("K26").Value = -4.347826086956
$ws.Range("L26").Value = -8.333333333333
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = 32.258064516129
$ws.Range("L27").Value = -2.380952380952
$ws.Range("N28").Value = -63.157894736842
$ws.Range("N29").Value = -66.666666666666

$excel.CutCopyMode = 0
